$d = $word.ActiveDocument

# The title paragraph is the first paragraph in the document.
$titlePara = $d.Paragraphs(1)

# Set the font size (in points) for the whole paragraph, including the
# paragraph mark, so both the run and the paragraph mark's run properties
# get sz/szCs = 28 (half-points) = 14pt.
$titlePara.Range.Font.Size = 14
$titlePara.Range.Font.SizeBi = 14
